# Insert a new data row at row 397 (pushing existing rows 397-457 down to 398-458)
# and populate it with the new reading for "Pepino ensalada".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(397).Insert()

$ws.Range("A397").Value = 9
$ws.Range("B397").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C397").Value = "Metropolitana"
$ws.Range("D397").Value = 45218
$ws.Range("E397").Value = 13
$ws.Range("F397").Value = 100112043
$ws.Range("G397").Value = "Pepino ensalada"
$ws.Range("H397").Value = "Sin especificar"
$ws.Range("I397").Value = "Primera"
$ws.Range("J397").Value = 70
$ws.Range("K397").Value = 12000
$ws.Range("L397").Value = 13000
$ws.Range("M397").Value = 12500
$ws.Range("N397").Value = "`$/caja 60 unidades"
$ws.Range("O397").Value = "Región de Arica y Parinacota"
$ws.Range("P397").Value = 208
$ws.Range("Q397").Value = 60
$ws.Range("R397").Value = "Hortaliza"
